$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.975.96'
$ws.Range("E2").Value = '  -1.36%  '

$ws.Range("D3").Value = '2.172.62'
$ws.Range("E3").Value = '  -2.41%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '66.28'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.62%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  -0.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.95'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0928'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.43%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.63'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -15.76%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.34%  '

$ws.Range("D15").Value = '2.483.62'
$ws.Range("E15").Value = '  -2.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.854'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.14%  '

$ws.Range("E17").Value = '  -4.44%  '

$ws.Range("D18").Value = '2.167.44'
$ws.Range("E18").Value = '  -2.62%  '

$ws.Range("D19").Value = '40.944.57'
$ws.Range("E19").Value = '  -1.47%  '

$ws.Range("D20").Value = '0.0₃0937'
$ws.Range("E20").Value = '  -2.97%  '

$ws.Range("E21").Value = '  -1.96%  '

$ws.Range("E22").Value = '  -2.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.13%  '

$ws.Range("E24").Value = '  -8.02%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.09%  '

$ws.Range("E27").Value = '  -4.85%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.55%  '

$ws.Range("E29").Value = '  -5.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.09%  '

$ws.Range("E33").Value = '  -0.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0753'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.47%  '

$ws.Range("E36").Value = '  -2.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.56'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0305'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.20%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.66%  '

$ws.Range("E42").Value = '  -8.87%  '

$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.61%  '

$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -13.37%  '

$ws.Range("E46").Value = '  -8.06%  '

$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("E48").Value = '  -4.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0988'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.77%  '

$ws.Range("E50").Value = '  -0.70%  '

$ws.Range("E51").Value = '  -3.87%  '
